# Scheduled-runner price refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across several crafting-sheet leve rows to the
# latest market-board figures.

$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param(
        [string]$SheetName,
        [int]$Row,
        $Values,
        $ClearCols = @()
    )

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($col in $Values.Keys) {
        $addr = "$col$Row"
        $ws.Range($addr).Value = $Values[$col]
    }

    foreach ($col in $ClearCols) {
        $addr = "$col$Row"
        $ws.Range($addr).ClearContents()
    }
}

# ---------------------------------------------------------------- ALC -----
Set-LeveRow "ALC" 15  @{ H=1746.1459; I=1746.1459; K=5238.4377; M=-5069.4377 }
Set-LeveRow "ALC" 112 @{ H=1161.5; I=273; J=1266.0294; K=819; L=3798.0882; M=289; N=-6014.0882 }
Set-LeveRow "ALC" 116 @{ H=4144.6895; I=2299.6667; J=4626; K=2299.6667; L=4626; M=1142.3333; N=-11510 }
Set-LeveRow "ALC" 118 @{ H=339.875; I=339.875; J=0; K=1019.625; L=0; M=637.375 } @("N")
Set-LeveRow "ALC" 132 @{ H=1208; I=854.4; J=5333.3335; K=2563.2; L=16000.0005; M=-33.19999999999982; N=-21060.0005 }

# ---------------------------------------------------------------- ARM -----
Set-LeveRow "ARM" 2   @{ H=1269.1177; I=727.7778; J=1878.125; K=727.7778; L=1878.125; M=-614.7778; N=-2104.125 }
Set-LeveRow "ARM" 32  @{ H=9226.946; I=3862.6743; J=26970.309; K=3862.6743; L=26970.309; M=-3575.6743; N=-27544.309 }
Set-LeveRow "ARM" 45  @{ H=1327.5385; I=1178; J=1566.8; K=1178; L=1566.8; M=-801; N=-2320.8 }
Set-LeveRow "ARM" 116 @{ H=1269.1177; I=727.7778; J=1878.125; K=727.7778; L=1878.125; M=1566.2222; N=-6466.125 }
Set-LeveRow "ARM" 132 @{ H=1858.2963; I=1006.6818; J=5605.4; K=3020.0454; L=16816.2; M=-490.0454; N=-21876.2 }

# ---------------------------------------------------------------- BSM -----
Set-LeveRow "BSM" 3   @{ H=1269.1177; I=727.7778; J=1878.125; K=727.7778; L=1878.125; M=-613.7778; N=-2106.125 }
Set-LeveRow "BSM" 24  @{ H=0; I=0; J=0; K=0; L=0 } @("M", "N")
Set-LeveRow "BSM" 86  @{ H=1483.25; I=1437; J=1622; K=1437; L=1622; M=-314; N=-3868 }
Set-LeveRow "BSM" 89  @{ H=1483.25; I=1437; J=1622; K=7185; L=8110; M=-1569; N=-19342 }
Set-LeveRow "BSM" 134 @{ H=2138.5588; I=2012.5518; J=2869.4; K=6037.6554; L=8608.200000000001; M=-3502.6554; N=-13678.2 }

# ---------------------------------------------------------------- CRP -----
Set-LeveRow "CRP" 31  @{ H=1437.4; I=856.7406999999999; J=6663.3335; K=856.7406999999999; L=6663.3335; M=-561.7406999999999; N=-7253.3335 }
Set-LeveRow "CRP" 34  @{ H=1437.4; I=856.7406999999999; J=6663.3335; K=856.7406999999999; L=6663.3335; M=-654.7406999999999; N=-7067.3335 }
Set-LeveRow "CRP" 54  @{ H=64500; I=0; J=64500; K=0; L=64500; N=-65816 }
Set-LeveRow "CRP" 86  @{ I=33082.688; J=3383.2727; K=33082.688; L=3383.2727; M=-31959.688; N=-5629.2727 }
Set-LeveRow "CRP" 89  @{ I=33082.688; J=3383.2727; K=165413.44; L=16916.3635; M=-159797.44; N=-28148.3635 }
Set-LeveRow "CRP" 99  @{ H=4008.6667; I=4012; J=4007; K=4012; L=4007; M=-2514; N=-7003 }
Set-LeveRow "CRP" 126 @{ H=4008.6667; I=4012; J=4007; K=12036; L=12021; M=-9566; N=-16961 }
Set-LeveRow "CRP" 132 @{ H=3544.7; I=3365.2856; J=3963.3333; K=10095.8568; L=11889.9999; M=-7565.856800000001; N=-16949.9999 }

# ---------------------------------------------------------------- GSM -----
Set-LeveRow "GSM" 138 @{ H=39143; I=39000; J=39214.5; K=39000; L=39214.5; M=-33860; N=-49494.5 }

# ---------------------------------------------------------------- LTW -----
Set-LeveRow "LTW" 30  @{ H=468.5; I=468.5; K=468.5; M=-360.5 }
Set-LeveRow "LTW" 130 @{ H=27000; J=27000; L=27000; N=-37040 }
Set-LeveRow "LTW" 132 @{ H=3287.1333; I=2250.4; J=5360.6; K=6751.200000000001; L=16081.8; M=-4221.200000000001; N=-21141.8 }

# ---------------------------------------------------------------- WVR -----
Set-LeveRow "WVR" 96  @{ H=3000; I=0; J=3000; K=0; L=3000; N=-5746 } @("M")
